$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-22 06:40:56'
$ws.Range('N2').Value = '0.6 °C 6:22 TU'
$ws.Range('O2').Value = '1.5 °C'
$ws.Range('E3').Value = '2026-02-22 06:40:59'
$ws.Range('E4').Value = '2026-02-22 06:41:01'
$ws.Range('O4').Value = '6.1 °C'
$ws.Range('E5').Value = '2026-02-22 06:41:04'
$ws.Range('O5').Value = '4.3 °C'
$ws.Range('E6').Value = '2026-02-22 06:41:07'
$ws.Range('H6').Value = '76%'
$ws.Range('E7').Value = '2026-02-22 06:41:09'
$ws.Range('E8').Value = '2026-02-22 06:41:12'
$ws.Range('H8').Value = '49%'
$ws.Range('O8').Value = '12.5 °C'
$ws.Range('E9').Value = '2026-02-22 06:41:15'
$ws.Range('H9').Value = '95%'
$ws.Range('E10').Value = '2026-02-22 06:41:17'
$ws.Range('N10').Value = '2.7 °C 6:15 TU'
$ws.Range('E11').Value = '2026-02-22 06:41:20'
$ws.Range('N11').Value = '0.7 °C 6:28 TU'
$ws.Range('O11').Value = '1.4 °C'
$ws.Range('E12').Value = '2026-02-22 06:41:23'
$ws.Range('N12').Value = '2.1 °C 6:04 TU'
$ws.Range('O12').Value = '4.2 °C'
$ws.Range('E13').Value = '2026-02-22 06:41:25'
$ws.Range('H13').Value = '92%'
$ws.Range('N13').Value = '-4.2 °C 6:04 TU'
$ws.Range('O13').Value = '-2.7 °C'
$ws.Range('E14').Value = '2026-02-22 06:41:28'
$ws.Range('H14').Value = '85%'
$ws.Range('N14').Value = '6.6 °C 6:05 TU'
$ws.Range('O14').Value = '7.4 °C'
$ws.Range('E15').Value = '2026-02-22 06:41:31'
$ws.Range('H15').Value = '88%'
$ws.Range('N15').Value = '2.0 °C 6:11 TU'
$ws.Range('O15').Value = '4.1 °C'
$ws.Range('E16').Value = '2026-02-22 06:41:33'
$ws.Range('L16').Value = '17.3 km/h - 218º 6:23 TU'
$ws.Range('M16').Value = '5.7 °C 6:19 TU'
$ws.Range('O16').Value = '4.2 °C'
$ws.Range('E17').Value = '2026-02-22 06:41:36'
$ws.Range('E18').Value = '2026-02-22 06:41:38'
$ws.Range('L18').Value = '5.8 km/h - 273º 6:15 TU'
$ws.Range('N18').Value = '-0.2 °C 6:27 TU'
$ws.Range('O18').Value = '1.1 °C'
$ws.Range('E19').Value = '2026-02-22 06:41:41'
$ws.Range('L19').Value = '7.6 km/h - 98º 6:20 TU'
$ws.Range('E20').Value = '2026-02-22 06:41:44'
$ws.Range('E21').Value = '2026-02-22 06:41:47'
$ws.Range('H21').Value = '79%'
$ws.Range('N21').Value = '0.4 °C 6:16 TU'
$ws.Range('O21').Value = '2.3 °C'
$ws.Range('E22').Value = '2026-02-22 06:41:49'
$ws.Range('L22').Value = '18.0 km/h - 315º 6:17 TU'
$ws.Range('E23').Value = '2026-02-22 06:41:52'
$ws.Range('M23').Value = '4.9 °C 6:15 TU'
$ws.Range('O23').Value = '3.9 °C'
$ws.Range('E24').Value = '2026-02-22 06:41:54'
$ws.Range('J24').Value = '1031.3 hPa'
$ws.Range('N24').Value = '-0.3 °C 6:27 TU'
$ws.Range('E25').Value = '2026-02-22 06:41:57'
$ws.Range('H25').Value = '27%'
$ws.Range('M25').Value = '5.6 °C 6:28 TU'
$ws.Range('O25').Value = '4.5 °C'
$ws.Range('E26').Value = '2026-02-22 06:42:00'
$ws.Range('E27').Value = '2026-02-22 06:42:03'
$ws.Range('H27').Value = '29%'
$ws.Range('O27').Value = '4.1 °C'
$ws.Range('E28').Value = '2026-02-22 06:42:06'
$ws.Range('O28').Value = '1.9 °C'
$ws.Range('E29').Value = '2026-02-22 06:42:08'
$ws.Range('N29').Value = '1.8 °C 6:27 TU'
$ws.Range('O29').Value = '4.1 °C'
$ws.Range('E30').Value = '2026-02-22 06:42:11'
$ws.Range('O30').Value = '7.4 °C'
$ws.Range('E31').Value = '2026-02-22 06:42:14'
$ws.Range('H31').Value = '68%'
$ws.Range('O31').Value = '11.7 °C'
$ws.Range('E32').Value = '2026-02-22 06:42:16'
$ws.Range('N32').Value = '-4.0 °C 6:13 TU'
$ws.Range('O32').Value = '-3.4 °C'
$ws.Range('E33').Value = '2026-02-22 06:42:19'
$ws.Range('N33').Value = '-0.5 °C 6:20 TU'
$ws.Range('O33').Value = '1.0 °C'
$ws.Range('E34').Value = '2026-02-22 06:42:21'
$ws.Range('O34').Value = '2.4 °C'
$ws.Range('E35').Value = '2026-02-22 06:42:24'
$ws.Range('H35').Value = '42%'
$ws.Range('E36').Value = '2026-02-22 06:42:26'
$ws.Range('E37').Value = '2026-02-22 06:42:29'
$ws.Range('J37').Value = '1033.7 hPa'
$ws.Range('L37').Value = '10.1 km/h - 31º 6:29 TU'
$ws.Range('M37').Value = '2.3 °C 6:29 TU'
$ws.Range('E38').Value = '2026-02-22 06:42:31'
$ws.Range('N38').Value = '2.6 °C 6:25 TU'
$ws.Range('O38').Value = '4.9 °C'
$ws.Range('E39').Value = '2026-02-22 06:42:34'
$ws.Range('L39').Value = '18.4 km/h - 273º 6:08 TU'
$ws.Range('E40').Value = '2026-02-22 06:42:37'
$ws.Range('H40').Value = '72%'
$ws.Range('N40').Value = '1.5 °C 6:20 TU'
$ws.Range('O40').Value = '4.4 °C'
$ws.Range('E41').Value = '2026-02-22 06:42:39'
$ws.Range('J41').Value = '1028.2 hPa'
$ws.Range('N41').Value = '3.6 °C 6:23 TU'
$ws.Range('O41').Value = '5.0 °C'
$ws.Range('E42').Value = '2026-02-22 06:42:42'
$ws.Range('O42').Value = '4.6 °C'
$ws.Range('E43').Value = '2026-02-22 06:42:45'
$ws.Range('K43').Value = '-0.1 MJ/m2'
$ws.Range('N43').Value = '-0.1 °C 6:28 TU'
$ws.Range('O43').Value = '1.7 °C'
$ws.Range('E44').Value = '2026-02-22 06:42:47'
$ws.Range('N44').Value = '-2.1 °C 6:24 TU'
$ws.Range('O44').Value = '-0.2 °C'
$ws.Range('E45').Value = '2026-02-22 06:42:50'
$ws.Range('H45').Value = '63%'
$ws.Range('J45').Value = '1031.7 hPa'
$ws.Range('N45').Value = '0.7 °C 6:09 TU'
$ws.Range('O45').Value = '4.1 °C'
$ws.Range('E46').Value = '2026-02-22 06:42:53'
$ws.Range('L46').Value = '5.4 km/h - 293º 6:14 TU'
$ws.Range('N46').Value = '-0.5 °C 6:05 TU'
$ws.Range('O46').Value = '1.0 °C'
